$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2 and C2 become completely empty (the TRUE flag in C2 is removed)
$ws.Range("B2").Clear()
$ws.Range("C2").Clear()

# Row 3: the TRUE flag moves from B3 to D3; B3 and C3 become completely empty
$ws.Range("B3").Clear()
$ws.Range("C3").Clear()
$ws.Range("D3").Value = $true

# Update the active selection to C2
$ws.Range("C2").Select()
